$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New handback file: f4a10675-5af9-4150-8d95-9d354567f2a2.md
# Adds one new row (row 3) to each of the three report sheets:
#   Sheet 1 "Overview", Sheet 2 "zh-cn", Sheet 3 "de-de"
# ---------------------------------------------------------------------

$newFile      = "f4a10675-5af9-4150-8d95-9d354567f2a2.md"
$newFileE2e   = "e2e\f4a10675-5af9-4150-8d95-9d354567f2a2.md"
$zhXlf        = "f4a10675-5af9-4150-8d95-9d354567f2a2.04ad14d99f32147bd293390af17307aa5afc71f8.zh-cn.xlf"
$deXlf        = "f4a10675-5af9-4150-8d95-9d354567f2a2.04ad14d99f32147bd293390af17307aa5afc71f8.de-de.xlf"

$srcUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/56e4c8eb81fe7d53ce3ec971a445a151ed221231/e2e/$newFile"
$zhUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/1f7f65c7ec67ead98b48c173dbd4f4c847504146/e2e/$newFile"
$deUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/bc339ed7faf7e07c33af5b3fccdd58c06e959461/e2e/$newFile"

# ============================= Sheet 1: Overview =======================
$ws1 = $wb.Worksheets.Item("Overview")
$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Range("A3").Value = $newFile
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2017-02-21 09:32:48"
$ws1.Range("G3").NumberFormat = $ws1.Range("G2").NumberFormat

$ws1.Hyperlinks.Add($ws1.Range("B3"), $srcUrl, "", "", $newFileE2e) | Out-Null

# ============================= Sheet 2: zh-cn ===========================
$ws2 = $wb.Worksheets.Item("zh-cn")
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "'True"
$ws2.Range("G3").Value = $zhXlf
$ws2.Range("H3").Value = "2017-02-21 09:32:32"
$ws2.Range("H3").NumberFormat = $ws2.Range("H2").NumberFormat
$ws2.Range("I3").Value = "'"
$ws2.Range("K3").Value = $zhXlf
$ws2.Range("L3").Value = "2017-02-21 09:33:31"
$ws2.Range("L3").NumberFormat = $ws2.Range("L2").NumberFormat
$ws2.Range("M3").Value = "'"
$ws2.Range("N3").Value = "'"
$ws2.Range("O3").Value = "'True"
$ws2.Range("P3").Value = "'"
$ws2.Range("Q3").Value = "'False"
$ws2.Range("R3").Value = "'"

$ws2.Hyperlinks.Add($ws2.Range("A3"), $srcUrl, "", "", $newFile) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("J3"), $zhUrl, "", "", $newFile) | Out-Null

# ============================= Sheet 3: de-de ===========================
$ws3 = $wb.Worksheets.Item("de-de")
$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "'True"
$ws3.Range("G3").Value = $deXlf
$ws3.Range("H3").Value = "2017-02-21 09:32:48"
$ws3.Range("H3").NumberFormat = $ws3.Range("H2").NumberFormat
$ws3.Range("I3").Value = "'"
$ws3.Range("K3").Value = $deXlf
$ws3.Range("L3").Value = "2017-02-21 09:33:56"
$ws3.Range("L3").NumberFormat = $ws3.Range("L2").NumberFormat
$ws3.Range("M3").Value = "'"
$ws3.Range("N3").Value = "'"
$ws3.Range("O3").Value = "'True"
$ws3.Range("P3").Value = "'"
$ws3.Range("Q3").Value = "'False"
$ws3.Range("R3").Value = "'"

$ws3.Hyperlinks.Add($ws3.Range("A3"), $srcUrl, "", "", $newFile) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("J3"), $deUrl, "", "", $newFile) | Out-Null

Write-Host "Generate Report for Handback: added row for $newFile to Overview, zh-cn, de-de"
